$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes existing rows 5-6 down to 6-7)
$ws.Rows.Item(5).Insert()

# New row 5 data: "Primera" quality, different price/date/unit info
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44427
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104005
$ws.Range("J5").Value = "Pera asiática"
$ws.Range("K5").Value = "Hosui"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("Q5").Value = "$/caja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 467
$ws.Range("T5").Value = 15

# Append a new row 8 at the end with fresh data (Especial quality)
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44418
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = "Frutos de pepita"
$ws.Range("I8").Value = 100104005
$ws.Range("J8").Value = "Pera asiática"
$ws.Range("K8").Value = "Hosui"
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("Q8").Value = "$/caja 15 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 533
$ws.Range("T8").Value = 15
